# The formula in column I (rows 2-31) of the "Data PGK B.Makanan (p)" sheet
# was updated so that the price-relatif factor is raised to the power of
# the economic-of-scale factor instead of being multiplied by it:
#   D{row} * 'Data PGK B.Makanan (β)'!$C$3 * 'Data PGK B.Makanan (β)'!$C$8
# becomes
#   D{row} * 'Data PGK B.Makanan (β)'!$C$3 ^ 'Data PGK B.Makanan (β)'!$C$8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data PGK B.Makanan (p)")

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("I$row")
    $oldFormula = $cell.Formula
    $lastStar = $oldFormula.LastIndexOf("*")
    $newFormula = $oldFormula.Substring(0, $lastStar) + "^" + $oldFormula.Substring($lastStar + 1)
    $cell.Formula = $newFormula
}
